# "fix: Reduce zoom level of excel file"
#
# The sheet was saved zoomed in to 120%. Bring the view back down to the
# standard 100% zoom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the sheet's zoom level back to 100%.
$excel.ActiveWindow.Zoom = 100

# Row 4 held an oversized 13.5pt font in Q4 which forced the row taller than
# its neighbours; now that the font is back to the normal 11pt size used
# elsewhere on the sheet, the row can shrink back down to match.
$q4 = $ws.Range("Q4")
$q4.Font.Name = "Arial"
$q4.Font.Size = 11
$q4.Font.Bold = $false
$q4.Font.Color = 0

$ws.Rows.Item(4).RowHeight = 14.25
